$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 490
$ws.Range("F5").Value = 920
$ws.Range("F6").Value = 154
$ws.Range("F7").Value = 938
$ws.Range("F8").Value = 731
$ws.Range("F9").Value = 179
$ws.Range("F10").Value = 47
$ws.Range("F12").Value = 768
$ws.Range("F13").Value = 255
$ws.Range("F14").Value = 549
$ws.Range("F15").Value = 488
$ws.Range("F16").Value = 1287
$ws.Range("F19").Value = 1072
$ws.Range("F20").Value = 2789
$ws.Range("F21").Value = 1274
$ws.Range("F23").Value = 165
$ws.Range("F24").Value = 1239
$ws.Range("F25").Value = 54
$ws.Range("F26").Value = 967
$ws.Range("F27").Value = 317
$ws.Range("F28").Value = 689
$ws.Range("F29").Value = 8
$ws.Range("F31").Value = 1313

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 348
$ws.Range("F7").Value = 5

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 490
$ws.Range("F9").Value = 348
$ws.Range("F12").Value = 920
$ws.Range("F13").Value = 154
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 938
$ws.Range("F16").Value = 731
$ws.Range("F17").Value = 179
$ws.Range("F19").Value = 47
$ws.Range("F25").Value = 768
$ws.Range("F26").Value = 255
$ws.Range("F27").Value = 549
$ws.Range("F28").Value = 488
$ws.Range("F29").Value = 1287
$ws.Range("F32").Value = 1072
$ws.Range("F33").Value = 2789
$ws.Range("F34").Value = 1274
$ws.Range("F36").Value = 165
$ws.Range("F37").Value = 1239
$ws.Range("F38").Value = 54
$ws.Range("F41").Value = 967
$ws.Range("F42").Value = 317
$ws.Range("F43").Value = 689
$ws.Range("F44").Value = 8
$ws.Range("F46").Value = 1313
